$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 257, shifting existing rows 257-274 down to 258-275
$ws.Rows.Item(257).Insert()

# Populate the newly inserted row 257 with the new record
$ws.Range("A257").Value = 5
$ws.Range("B257").Value = "Macroferia Regional de Talca"
$ws.Range("C257").Value = "Maule"
$ws.Range("D257").Value = 44931
$ws.Range("E257").Value = 7
$ws.Range("F257").Value = 100112021
$ws.Range("G257").Value = "Ají"
$ws.Range("H257").Value = "Americana (o)"
$ws.Range("I257").Value = "Primera"
$ws.Range("J257").Value = 150
$ws.Range("K257").Value = 12000
$ws.Range("L257").Value = 12000
$ws.Range("M257").Value = 12000
$ws.Range("N257").Value = "$/caja 15 kilos"
$ws.Range("O257").Value = "Región del Maule"
$ws.Range("P257").Value = 800
$ws.Range("Q257").Value = 15
$ws.Range("R257").Value = "Hortaliza"
